$p = $ppt.ActivePresentation

# 1. Delete the second slide (the duplicated "sticker strip" picture slide
#    that was only reintroduced by the previous hotfix).
$p.Slides.Item(2).Delete()

# 2. Re-merge the two runs of the second paragraph in the "Differences..."
#    rectangle on slide 1 back into a single run (undoing the accidental
#    run-split), same net text either way.
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(18)
$tr = $shape.TextFrame.TextRange
$para = $tr.Paragraphs(2, 1)
# Force the underlying runs to coalesce: assign a scratch value first so the
# subsequent assignment of the real text is recognised as a genuine edit
# (identical before/after text is otherwise treated as a no-op and the
# original run split would survive untouched).
$para.Text = "_tmp_"
$tr.Paragraphs(2, 1).Text = "Differences cannot be explained by sharing distributional info."
